$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-03-28 Friday" "2025-03-29 Saturday"

Replace-Text "37×21=" "56×19="
Replace-Text "71×95=" "56×20="
Replace-Text "98×98=" "27×19="
Replace-Text "89×25=" "88×32="
Replace-Text "93×34=" "80×30="
Replace-Text "32×76=" "21×60="
Replace-Text "27×20=" "52×17="
Replace-Text "57×74=" "74×72="
Replace-Text "91×12=" "81×94="
Replace-Text "61×19=" "50×54="
Replace-Text "60×62=" "88×18="
Replace-Text "88×46=" "90×72="
Replace-Text "36×78=" "67×66="
Replace-Text "73×87=" "31×27="
Replace-Text "71×94=" "74×11="
Replace-Text "30×31=" "82×81="
Replace-Text "49×86=" "26×42="
Replace-Text "93×18=" "19×85="
Replace-Text "89×34=" "79×11="
Replace-Text "93×65=" "47×26="
Replace-Text "56×67=" "25×29="
Replace-Text "51×98=" "94×30="
Replace-Text "35×91=" "30×46="
Replace-Text "75×59=" "39×39="
Replace-Text "13×68=" "62×79="
